# Append a new listing ("PHP・Laravel" job) as the new row 10 of the
# "ランサーズ" sheet, pushing the previous rows 10-18 down to 11-19, and
# refresh the "取得日時" scrape timestamp on every data row to the latest
# run (2025-09-22 18:24:26).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# The emulated Hyperlinks collection does not re-target itself when rows are
# inserted/shifted, so clear it first and rebuild it from scratch once all
# the row data is back in its final place. This keeps every hyperlink's
# worksheet ref in sync with the relationship it points to.
$ws.Hyperlinks.Delete()

# Inserting a row at 10 shifts the old rows 10-18 (and their formatting) down
# to 11-19, and grows the sheet dimension to A1:H19 automatically.
$ws.Rows.Item(10).Insert()

# Refresh the scrape timestamp in column A for every data row (2-19).
$newTimestamp = "2025-09-22 18:24:26"
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Populate the brand-new row 10 with the newly scraped listing.
$ws.Cells.Item(10, 2).Value = "【急募】PHP・Lalavelでの既存プログラム改修依頼"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5396563"
$ws.Cells.Item(10, 7).Value = 33
$ws.Cells.Item(10, 8).Value = "○PHP"

# Rebuild the column-F hyperlinks for every data row, in order, so the
# worksheet's <hyperlinks> entries (and the rels they point at) line up with
# the URL typed into each F cell above.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5398198")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5398112")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5398193")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5398203")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5397930")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5397812")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5398081")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5396563")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5398382")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5398071")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5398062")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5398293")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5397887")
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5398497")
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5397980")
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5397962")
$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.lancers.jp/work/detail/5397817")
